# DB Validation added for drivers
# - Shorten the credit-card segmentation label on the Filter sheet
# - Make "up_sell_Filter" (first sheet) the active/selected tab instead of
#   "up_Sell_Report_EtoE" (second sheet)

$wb = $excel.ActiveWorkbook

# up_sell_Filter is the first worksheet in the workbook
$wsFilter = $wb.Worksheets.Item(1)

# Update the Credit Cards segmentation text (was a long descriptive label,
# now just "Credit Cards")
$wsFilter.Range("C2").Value = "Credit Cards"

# Make this sheet the active tab (moves tabSelected from the report sheet
# to the filter sheet, and resets the workbook's active tab)
$wsFilter.Activate()
